$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1862
$ws.Cells.Item(5, 6).Value = 59
$ws.Cells.Item(8, 6).Value = 660
$ws.Cells.Item(12, 6).Value = 1508
$ws.Cells.Item(13, 6).Value = 1228
$ws.Cells.Item(14, 6).Value = 1465
$ws.Cells.Item(15, 6).Value = 36
$ws.Cells.Item(16, 6).Value = 1313
$ws.Cells.Item(17, 6).Value = 307
$ws.Cells.Item(18, 6).Value = 1613
$ws.Cells.Item(20, 6).Value = 784
$ws.Cells.Item(21, 6).Value = 1038
$ws.Cells.Item(22, 6).Value = 342
$ws.Cells.Item(23, 6).Value = 49
$ws.Cells.Item(24, 6).Value = 106
$ws.Cells.Item(25, 6).Value = 1482
$ws.Cells.Item(27, 6).Value = 150
$ws.Cells.Item(29, 6).Value = 550
$ws.Cells.Item(30, 6).Value = 1112
$ws.Cells.Item(31, 6).Value = 996
$ws.Cells.Item(32, 6).Value = 39
$ws.Cells.Item(33, 6).Value = 559
$ws.Cells.Item(35, 6).Value = 1077
$ws.Cells.Item(36, 6).Value = 891
$ws.Cells.Item(37, 6).Value = 1084
$ws.Cells.Item(38, 6).Value = 39
$ws.Cells.Item(39, 6).Value = 119
$ws.Cells.Item(41, 6).Value = 845
$ws.Cells.Item(42, 6).Value = 1621
$ws.Cells.Item(44, 6).Value = 61
$ws.Cells.Item(45, 6).Value = 799
$ws.Cells.Item(47, 6).Value = 781

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 44
$ws.Cells.Item(6, 6).Value = 178
$ws.Cells.Item(7, 6).Value = 1466
$ws.Cells.Item(10, 6).Value = 2548
$ws.Cells.Item(14, 6).Value = 228
$ws.Cells.Item(19, 6).Value = 443
$ws.Cells.Item(23, 6).Value = 74388
$ws.Cells.Item(28, 6).Value = 239
$ws.Cells.Item(30, 6).Value = 169
$ws.Cells.Item(32, 6).Value = 31
$ws.Cells.Item(35, 6).Value = 174

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 2836
$ws.Cells.Item(6, 6).Value = 4566
$ws.Cells.Item(9, 6).Value = 550
$ws.Cells.Item(10, 6).Value = 679
$ws.Cells.Item(11, 6).Value = 441
$ws.Cells.Item(12, 6).Value = 267
$ws.Cells.Item(13, 6).Value = 891
$ws.Cells.Item(14, 6).Value = 231
$ws.Cells.Item(15, 6).Value = 537

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1862
$ws.Cells.Item(4, 6).Value = 2836
$ws.Cells.Item(5, 6).Value = 4566
$ws.Cells.Item(6, 6).Value = 679
$ws.Cells.Item(7, 6).Value = 441
$ws.Cells.Item(8, 6).Value = 267
$ws.Cells.Item(9, 6).Value = 267
$ws.Cells.Item(10, 6).Value = 891
$ws.Cells.Item(11, 6).Value = 891
$ws.Cells.Item(16, 6).Value = 2548
$ws.Cells.Item(18, 6).Value = 1508
$ws.Cells.Item(19, 6).Value = 1228
$ws.Cells.Item(20, 6).Value = 1465
$ws.Cells.Item(21, 6).Value = 1313
$ws.Cells.Item(22, 6).Value = 228
$ws.Cells.Item(23, 6).Value = 307
$ws.Cells.Item(25, 6).Value = 1613
$ws.Cells.Item(26, 6).Value = 784
$ws.Cells.Item(27, 6).Value = 1038
$ws.Cells.Item(28, 6).Value = 342
$ws.Cells.Item(29, 6).Value = 537
$ws.Cells.Item(30, 6).Value = 537
$ws.Cells.Item(31, 6).Value = 1482
$ws.Cells.Item(33, 6).Value = 150
$ws.Cells.Item(35, 6).Value = 550
$ws.Cells.Item(36, 6).Value = 1112
$ws.Cells.Item(38, 6).Value = 996
$ws.Cells.Item(39, 6).Value = 39
$ws.Cells.Item(40, 6).Value = 1077
$ws.Cells.Item(41, 6).Value = 891
$ws.Cells.Item(42, 6).Value = 1084
$ws.Cells.Item(44, 6).Value = 119
$ws.Cells.Item(45, 6).Value = 845
$ws.Cells.Item(47, 6).Value = 1621
$ws.Cells.Item(49, 6).Value = 61
$ws.Cells.Item(50, 6).Value = 799
$ws.Cells.Item(52, 6).Value = 781
